# "can read new excel format"
# The Percent column (E) was previously storing raw numbers (50, 75, 85, ...)
# formatted as General. The sheet now reads the new source data as true
# fractional percentages (0.5, 0.7, 0.85, ...) displayed with a "0%"
# number format for the rows that changed (Sprint 1 / Sprint 2 / Sprint 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sprint 1 / Prepare): 50 -> 50%
$ws.Range("E2").Value = 0.5
$ws.Range("E2").NumberFormat = "0%"

# Row 3 (Sprint 2 / Doing Work 1): 75 -> 70%
$ws.Range("E3").Value = 0.7
$ws.Range("E3").NumberFormat = "0%"

# Row 8 (Sprint 7 / Doing Work 6): 85 -> 85%
$ws.Range("E8").Value = 0.85
$ws.Range("E8").NumberFormat = "0%"

# Bring the view back to the top of the sheet and move the active
# selection to E9, matching the refreshed view state of the workbook.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E9").Select() | Out-Null
